# Auto update: 2025-05-20 20:07:55
# Apply row-content updates to columns A (Company Name), B (Company Number)
# and H (Category) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Company Name (A), Company Number (B), [optional] Category (H)
$updates = @(
    @{ Row = 2;  A = "DGPI LTD"; B = "SC849118"; H = "GP" },
    @{ Row = 3;  A = "DAVIDSON CAPITAL HOLDINGS LTD"; B = "SC849117"; H = "Capital" },
    @{ Row = 4;  A = "ST GEORGE CAPITAL (LAND) LIMITED"; B = "16462880" },
    @{ Row = 5;  A = "AFROSCOT VENTURES LTD"; B = "16462878"; H = "Ventures" },
    @{ Row = 6;  A = "SAMVIV PARTNERS LTD"; B = "16460672"; H = "Partners" },
    @{ Row = 7;  A = "T GILPIN PHYSIO CONSULTANCY LTD"; B = "16460503"; H = "LP" },
    @{ Row = 8;  A = "4D CAPITAL PROPCO (44) LIMITED"; B = "16461269" },
    @{ Row = 9;  A = "KNOTT INVESTMENTS LIMITED"; B = "16458684" },
    @{ Row = 10; A = "THE REEL MED LLP"; B = "OC456780"; H = "LP" },
    @{ Row = 11; A = "PONGPONG MALATANG LTD"; B = "16458077"; H = "GP" },
    @{ Row = 12; A = "KC INVESTMENTS & TRADING LIMITED"; B = "16456642" },
    @{ Row = 13; A = "JJOHN INVESTMENTS LIMITED"; B = "16456276"; H = "Investments" },
    @{ Row = 14; A = "ECHO VENTURES GROUP LIMITED"; B = "16455744"; H = "Ventures" },
    @{ Row = 16; A = "MUSICROOTS LTD"; B = "16455514"; H = "SIC" },
    @{ Row = 17; A = "DAVISON FAMILY CAPITAL LTD"; B = "16455115"; H = "Capital" },
    @{ Row = 19; A = "TALKSGPT AI LTD"; B = "16455313"; H = "GP" },
    @{ Row = 20; A = "SLAM DUNK INVESTMENTS LTD"; B = "16455167"; H = "Investments" },
    @{ Row = 21; A = "IX PARTNERS LLP"; B = "OC456771"; H = "Partners" },
    @{ Row = 23; A = "GROWTHFORGE MANAGEMENT LLP"; B = "OC456769"; H = "LP" },
    @{ Row = 24; A = "GOLDEN VENTURES LONDON LTD"; B = "16452104"; H = "Ventures" },
    @{ Row = 25; A = "ALDABBOUS UK INVESTMENTS LTD"; B = "16453476" },
    @{ Row = 26; A = "CAMBRIDGE SOCIAL INVESTMENTS LIMITED"; B = "16453466" },
    @{ Row = 27; A = "CAPITAL & CENTRIC (SYNCHRONICITY) LTD"; B = "16453716"; H = "Capital" },
    @{ Row = 28; A = "GULF TRADE AND INVESTMENT ADVANTAGES JOINT PARTNERSHIP LTD"; B = "16453733"; H = "Partners" },
    @{ Row = 29; A = "FROST CAPITAL LTD"; B = "16450073"; H = "Capital" },
    @{ Row = 31; A = "SYNERGY FUNDING LTD"; B = "16449538"; H = "Fund" },
    @{ Row = 32; A = "ASSET CAPITAL 44 OPCO LIMITED"; B = "16449512" }
)

# Company Number values that are purely digits would be auto-converted to a
# numeric cell type by plain `.Value` assignment (standard Excel behaviour).
# The source data keeps these as text cells (matching the original file), so
# for those rows we briefly mark the cell as Text before writing the value,
# then clear the formatting back off again so no visible/number-format
# change is left behind - only the underlying text cell type sticks.
$bCellsNeedingTextCoercion = @()
foreach ($u in $updates) {
    if ($u.B -match '^[0-9]+$') {
        $bCellsNeedingTextCoercion += $u.Row
    }
}

foreach ($r in $bCellsNeedingTextCoercion) {
    $ws.Cells.Item($r, 2).NumberFormat = "@"
}

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A
    $ws.Cells.Item($r, 2).Value = $u.B
    if ($u.ContainsKey("H")) {
        $ws.Cells.Item($r, 8).Value = $u.H
    }
}

foreach ($r in $bCellsNeedingTextCoercion) {
    $ws.Cells.Item($r, 2).ClearFormats()
}
